{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// practice table with its new value. Every original equation string is\n// unique within the document, so a plain body-wide search + insertText\n// \"replace\" for each (old, new) pair reproduces the diff while leaving\n// every other run property (font, size, etc.) untouched.\nconst replacements = [\n  [\"206\u00d74=824\", \"987\u00d74=3948\"],\n  [\"829\u00d73=2487\", \"232\u00d77=1624\"],\n  [\"499\u00d75=2495\", \"645\u00d77=4515\"],\n  [\"536\u00d72=1072\", \"498\u00d72=996\"],\n  [\"907\u00d72=1814\", \"306\u00d75=1530\"],\n  [\"751\u00d78=6008\", \"831\u00d76=4986\"],\n  [\"379\u00d78=3032\", \"386\u00d76=2316\"],\n  [\"898\u00d75=4490\", \"350\u00d74=1400\"],\n  [\"282\u00d75=1410\", \"872\u00d73=2616\"],\n  [\"872\u00d78=6976\", \"742\u00d78=5936\"],\n  [\"343\u00d75=1715\", \"892\u00d73=2676\"],\n  [\"956\u00d73=2868\", \"891\u00d73=2673\"],\n  [\"264\u00d73=792\", \"676\u00d75=3380\"],\n  [\"638\u00d76=3828\", \"121\u00d75=605\"],\n  [\"222\u00d73=666\", \"832\u00d73=2496\"],\n  [\"155\u00d75=775\", \"566\u00d79=5094\"],\n  [\"307\u00d78=2456\", \"894\u00d72=1788\"],\n  [\"215\u00d75=1075\", \"857\u00d77=5999\"],\n  [\"116\u00d74=464\", \"121\u00d76=726\"],\n  [\"138\u00d72=276\", \"784\u00d75=3920\"],\n  [\"932\u00d77=6524\", \"634\u00d75=3170\"],\n  [\"850\u00d72=1700\", \"998\u00d74=3992\"],\n  [\"400\u00d78=3200\", \"922\u00d77=6454\"],\n  [\"413\u00d74=1652\", \"233\u00d78=1864\"],\n  [\"958\u00d78=7664\", \"900\u00d74=3600\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each equation in the practice table is unique, so a straight\n# Find/Replace (wdReplaceAll) over the whole document body for each\n# old->new pair reproduces the diff while leaving every other run\n# property (font, size, etc.) untouched.\n$pairs = @(\n  @('206\u00d74=824', '987\u00d74=3948'),\n  @('829\u00d73=2487', '232\u00d77=1624'),\n  @('499\u00d75=2495', '645\u00d77=4515'),\n  @('536\u00d72=1072', '498\u00d72=996'),\n  @('907\u00d72=1814', '306\u00d75=1530'),\n  @('751\u00d78=6008', '831\u00d76=4986'),\n  @('379\u00d78=3032', '386\u00d76=2316'),\n  @('898\u00d75=4490', '350\u00d74=1400'),\n  @('282\u00d75=1410', '872\u00d73=2616'),\n  @('872\u00d78=6976', '742\u00d78=5936'),\n  @('343\u00d75=1715', '892\u00d73=2676'),\n  @('956\u00d73=2868', '891\u00d73=2673'),\n  @('264\u00d73=792', '676\u00d75=3380'),\n  @('638\u00d76=3828', '121\u00d75=605'),\n  @('222\u00d73=666', '832\u00d73=2496'),\n  @('155\u00d75=775', '566\u00d79=5094'),\n  @('307\u00d78=2456', '894\u00d72=1788'),\n  @('215\u00d75=1075', '857\u00d77=5999'),\n  @('116\u00d74=464', '121\u00d76=726'),\n  @('138\u00d72=276', '784\u00d75=3920'),\n  @('932\u00d77=6524', '634\u00d75=3170'),\n  @('850\u00d72=1700', '998\u00d74=3992'),\n  @('400\u00d78=3200', '922\u00d77=6454'),\n  @('413\u00d74=1652', '233\u00d78=1864'),\n  @('958\u00d78=7664', '900\u00d74=3600'),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $found = $find.Execute(\n    $oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2\n  )\n\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n\n"}
